$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<have>"
$ws.Range("C2").Value = 32

$ws.Range("C3").Value = 40

$ws.Range("B4").Value = "<no>"
$ws.Range("C4").Value = 44

$ws.Range("C5").Value = 31

$ws.Range("C6").Value = 38

$ws.Range("C7").Value = 36

$ws.Range("B8").Value = "<then>"
$ws.Range("C8").Value = 41

$ws.Range("C9").Value = 39

$ws.Range("B10").Value = "<we>"
$ws.Range("C10").Value = 36

$ws.Range("B11").Value = "<cope>"

$ws.Range("C12").Value = 38

$ws.Range("C13").Value = 39

$ws.Range("B14").Value = "<nome>"
$ws.Range("C14").Value = 41

$ws.Range("C15").Value = 38

$ws.Range("B16").Value = "<not>"
$ws.Range("C16").Value = 38

$ws.Range("B17").Value = "<enter>"
$ws.Range("C17").Value = 40

$ws.Range("B18").Value = "<we>"
$ws.Range("C18").Value = 28
